$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The question-order table in column F (and the "comment" column C) needs each
# entry shifted down by one row, i.e. a single blank row is effectively
# inserted above row 6 for columns C:F only (row 5 keeps its A/B cells).
# Cut/paste single cells bottom-up so we never clobber a cell we still need
# to move.
$ws.Range("F12").Cut($ws.Range("F13"))
$ws.Range("F11").Cut($ws.Range("F12"))
$ws.Range("F10").Cut($ws.Range("F11"))
$ws.Range("F9").Cut($ws.Range("F10"))
$ws.Range("C9").Cut($ws.Range("C10"))
$ws.Range("F8").Cut($ws.Range("F9"))
$ws.Range("F7").Cut($ws.Range("F8"))
$ws.Range("F6").Cut($ws.Range("F7"))
$ws.Range("C5:F5").Cut($ws.Range("C6:F6"))

# Drop the now-empty, but still formatted, tail of row 5 so it disappears
# completely instead of leaving a style-only cell behind.
$ws.Range("E5").NumberFormat = "General"
$ws.Range("E5").ClearContents()

# Row heights: rows 5 & 6 shrink slightly, and a trailing (otherwise empty)
# row at the very bottom of the sheet gets a height too.
$ws.Rows("5:5").RowHeight = 13.8
$ws.Rows("6:6").RowHeight = 13.8
$ws.Rows("1048576:1048576").RowHeight = 12.8

# Move the active selection.
[void]$ws.Range("A13").Select()
